# Insert a new data row at row 150 (this pushes the former rows 150-207
# down to 151-208, matching the diff's dimension change A1:T207 -> A1:T208
# and the shift of all subsequent rows' contents by one position).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(150).Insert()

# Populate the newly inserted row 150 with its values.
$ws.Cells.Item(150, 1).Value = 7
$ws.Cells.Item(150, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(150, 3).Value = "Ñuble"
$ws.Cells.Item(150, 4).Value = 45229
$ws.Cells.Item(150, 5).Value = 16
$ws.Cells.Item(150, 6).Value = "Fruta"
$ws.Cells.Item(150, 7).Value = 100108
$ws.Cells.Item(150, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(150, 9).Value = 100108002
$ws.Cells.Item(150, 10).Value = "Mango"
$ws.Cells.Item(150, 11).Value = "Sin especificar"
$ws.Cells.Item(150, 12).Value = "Primera"
$ws.Cells.Item(150, 13).Value = 60
$ws.Cells.Item(150, 14).Value = 10000
$ws.Cells.Item(150, 15).Value = 10000
$ws.Cells.Item(150, 16).Value = 10000
$ws.Cells.Item(150, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(150, 18).Value = "Brasil"
$ws.Cells.Item(150, 19).Value = 2500
$ws.Cells.Item(150, 20).Value = 4
